$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.71718733333334
$ws.Range("H2").Value = 110.151562
$ws.Range("I2").Value = 0.728743057485239
$ws.Range("J2").Value = 0.7328478088626956
$ws.Range("M2").Value = 0.5804443333333333
$ws.Range("N2").Value = 1.741333
$ws.Range("O2").Value = 0.002431273010151717
$ws.Range("P2").Value = 0.002435427107574628
$ws.Range("Q2").Value = 21.31228332357178
$ws.Range("R2").Value = 191.810549912146
$ws.Range("S2").Value = 0.001771773326999303
$ws.Range("T2").Value = 0.001784797419430879
$ws.Range("G3").Value = 36.71718733333334
$ws.Range("H3").Value = 110.151562
$ws.Range("I3").Value = 0.728743057485239
$ws.Range("J3").Value = 0.7328478088626956
$ws.Range("O3").Value = 0.0004752041289926495
$ws.Range("P3").Value = 0.00047601606752829
$ws.Range("Q3").Value = 4.165589381091555
$ws.Range("R3").Value = 37.490304429824
$ws.Range("S3").Value = 0.0003463017098917133
$ws.Range("T3").Value = 0.0003488473320715443
$ws.Range("G4").Value = 36.71718733333334
$ws.Range("H4").Value = 110.151562
$ws.Range("I4").Value = 0.728743057485239
$ws.Range("J4").Value = 0.7328478088626956
$ws.Range("M4").Value = 136.1000366666667
$ws.Range("N4").Value = 408.30011
$ws.Range("O4").Value = 0.5700742118164518
$ws.Range("P4").Value = 0.5710482463260632
$ws.Range("Q4").Value = 4997.210542363537
$ws.Range("R4").Value = 44974.89488127182
$ws.Range("S4").Value = 0.4154376241126089
$ws.Range("T4").Value = 0.4184914560749403
$ws.Range("G5").Value = 36.71718733333334
$ws.Range("H5").Value = 110.151562
$ws.Range("I5").Value = 0.728743057485239
$ws.Range("J5").Value = 0.7328478088626956
$ws.Range("M5").Value = 1.221658
$ws.Range("N5").Value = 2.443316
$ws.Range("O5").Value = 0.005117086949542552
$ws.Range("P5").Value = 0.003417220037046797
$ws.Range("Q5").Value = 44.85584564326533
$ws.Range("R5").Value = 269.135073859592
$ws.Range("S5").Value = 0.003729041589027455
$ws.Range("T5").Value = 0.002504302216551445
$ws.Range("G6").Value = 36.71718733333334
$ws.Range("H6").Value = 110.151562
$ws.Range("I6").Value = 0.728743057485239
$ws.Range("J6").Value = 0.7328478088626956
$ws.Range("M6").Value = 100.7253213333333
$ws.Range("N6").Value = 302.175964
$ws.Range("O6").Value = 0.4219022240948613
$ws.Range("P6").Value = 0.4226230904617871
$ws.Range("Q6").Value = 3698.350492606197
$ws.Range("R6").Value = 33285.15443345577
$ws.Range("S6").Value = 0.3074583167467117
$ws.Range("T6").Value = 0.3097184058197014
$ws.Range("I7").Value = 0.02522574977045663
$ws.Range("J7").Value = 0.0253678374789488
$ws.Range("M7").Value = 0.5804443333333333
$ws.Range("N7").Value = 1.741333
$ws.Range("O7").Value = 0.002431273010151717
$ws.Range("P7").Value = 0.002435427107574628
$ws.Range("Q7").Value = 0.7377337192243333
$ws.Range("R7").Value = 6.639603473018999
$ws.Range("S7").Value = [double]"6.133068457775209E-05"
$ws.Range("T7").Value = [double]"6.178151905677954E-05"
$ws.Range("I8").Value = 0.02522574977045663
$ws.Range("J8").Value = 0.0253678374789488
$ws.Range("O8").Value = 0.0004752041289926495
$ws.Range("P8").Value = 0.00047601606752829
$ws.Range("S8").Value = [double]"1.198738044785637E-05"
$ws.Range("T8").Value = [double]"1.207549823842598E-05"
$ws.Range("I9").Value = 0.02522574977045663
$ws.Range("J9").Value = 0.0253678374789488
$ws.Range("M9").Value = 136.1000366666667
$ws.Range("N9").Value = 408.30011
$ws.Range("O9").Value = 0.5700742118164518
$ws.Range("P9").Value = 0.5710482463260632
$ws.Range("Q9").Value = 172.9805607026367
$ws.Range("R9").Value = 1556.82504632373
$ws.Range("S9").Value = 0.01438054941787211
$ws.Range("T9").Value = 0.01448625910543829
$ws.Range("I10").Value = 0.02522574977045663
$ws.Range("J10").Value = 0.0253678374789488
$ws.Range("M10").Value = 1.221658
$ws.Range("N10").Value = 2.443316
$ws.Range("O10").Value = 0.005117086949542552
$ws.Range("P10").Value = 0.003417220037046797
$ws.Range("Q10").Value = 1.552704106498
$ws.Range("R10").Value = 9.316224638987999
$ws.Range("S10").Value = 0.0001290823549428297
$ws.Range("T10").Value = [double]"8.668748252961057E-05"
$ws.Range("I11").Value = 0.02522574977045663
$ws.Range("J11").Value = 0.0253678374789488
$ws.Range("M11").Value = 100.7253213333333
$ws.Range("N11").Value = 302.175964
$ws.Range("O11").Value = 0.4219022240948613
$ws.Range("P11").Value = 0.4226230904617871
$ws.Range("Q11").Value = 128.0199696335613
$ws.Range("R11").Value = 1152.179726702052
$ws.Range("S11").Value = 0.01064279993261609
$ws.Range("T11").Value = 0.01072103387368569
$ws.Range("G12").Value = 7.275657333333332
$ws.Range("H12").Value = 21.826972
$ws.Range("I12").Value = 0.1444033477339586
$ws.Range("J12").Value = 0.1452167206154317
$ws.Range("M12").Value = 0.5804443333333333
$ws.Range("N12").Value = 1.741333
$ws.Range("O12").Value = 0.002431273010151717
$ws.Range("P12").Value = 0.002435427107574628
$ws.Range("Q12").Value = 4.223114070408444
$ws.Range("R12").Value = 38.008026633676
$ws.Range("S12").Value = 0.0003510839619211267
$ws.Range("T12").Value = 0.0003536647378599138
$ws.Range("G13").Value = 7.275657333333332
$ws.Range("H13").Value = 21.826972
$ws.Range("I13").Value = 0.1444033477339586
$ws.Range("J13").Value = 0.1452167206154317
$ws.Range("O13").Value = 0.0004752041289926495
$ws.Range("P13").Value = 0.00047601606752829
$ws.Range("Q13").Value = 0.8254281749048887
$ws.Range("R13").Value = 7.428853574143999
$ws.Range("S13").Value = [double]"6.862106708353848E-05"
$ws.Range("T13").Value = [double]"6.912549228671218E-05"
$ws.Range("G14").Value = 7.275657333333332
$ws.Range("H14").Value = 21.826972
$ws.Range("I14").Value = 0.1444033477339586
$ws.Range("J14").Value = 0.1452167206154317
$ws.Range("M14").Value = 136.1000366666667
$ws.Range("N14").Value = 408.30011
$ws.Range("O14").Value = 0.5700742118164518
$ws.Range("P14").Value = 0.5710482463260632
$ws.Range("Q14").Value = 990.2172298407688
$ws.Range("R14").Value = 8911.95506856692
$ws.Range("S14").Value = 0.08232062464309346
$ws.Range("T14").Value = 0.08292575364466416
$ws.Range("G15").Value = 7.275657333333332
$ws.Range("H15").Value = 21.826972
$ws.Range("I15").Value = 0.1444033477339586
$ws.Range("J15").Value = 0.1452167206154317
$ws.Range("M15").Value = 1.221658
$ws.Range("N15").Value = 2.443316
$ws.Range("O15").Value = 0.005117086949542552
$ws.Range("P15").Value = 0.003417220037046797
$ws.Range("Q15").Value = 8.888364986525332
$ws.Range("R15").Value = 53.33018991915199
$ws.Range("S15").Value = 0.0007389244861596946
$ws.Range("T15").Value = 0.00049623748740128
$ws.Range("G16").Value = 7.275657333333332
$ws.Range("H16").Value = 21.826972
$ws.Range("I16").Value = 0.1444033477339586
$ws.Range("J16").Value = 0.1452167206154317
$ws.Range("M16").Value = 100.7253213333333
$ws.Range("N16").Value = 302.175964
$ws.Range("O16").Value = 0.4219022240948613
$ws.Range("P16").Value = 0.4226230904617871
$ws.Range("Q16").Value = 732.842922811223
$ws.Range("R16").Value = 6595.586305301008
$ws.Range("S16").Value = 0.06092409357570078
$ws.Range("T16").Value = 0.06137193925321967
$ws.Range("G17").Value = 0.8466215
$ws.Range("H17").Value = 1.693243
$ws.Range("I17").Value = 0.01680328982832053
$ws.Range("J17").Value = 0.0112652912032432
$ws.Range("M17").Value = 0.5804443333333333
$ws.Range("N17").Value = 1.741333
$ws.Range("O17").Value = 0.002431273010151717
$ws.Range("P17").Value = 0.002435427107574628
$ws.Range("Q17").Value = 0.4914166521531667
$ws.Range("R17").Value = 2.948499912919
$ws.Range("S17").Value = [double]"4.085338504135259E-05"
$ws.Range("T17").Value = [double]"2.743579557110048E-05"
$ws.Range("G18").Value = 0.8466215
$ws.Range("H18").Value = 1.693243
$ws.Range("I18").Value = 0.01680328982832053
$ws.Range("J18").Value = 0.0112652912032432
$ws.Range("O18").Value = 0.0004752041289926495
$ws.Range("P18").Value = 0.00047601606752829
$ws.Range("Q18").Value = 0.09604977358933332
$ws.Range("R18").Value = 0.576298641536
$ws.Range("S18").Value = [double]"7.984992707078104E-06"
$ws.Range("T18").Value = [double]"5.362459618128865E-06"
$ws.Range("G19").Value = 0.8466215
$ws.Range("H19").Value = 1.693243
$ws.Range("I19").Value = 0.01680328982832053
$ws.Range("J19").Value = 0.0112652912032432
$ws.Range("M19").Value = 136.1000366666667
$ws.Range("N19").Value = 408.30011
$ws.Range("O19").Value = 0.5700742118164518
$ws.Range("P19").Value = 0.5710482463260632
$ws.Range("Q19").Value = 115.2252171927883
$ws.Range("R19").Value = 691.3513031567301
$ws.Range("S19").Value = 0.009579122204803227
$ws.Range("T19").Value = 0.006433024785964453
$ws.Range("G20").Value = 0.8466215
$ws.Range("H20").Value = 1.693243
$ws.Range("I20").Value = 0.01680328982832053
$ws.Range("J20").Value = 0.0112652912032432
$ws.Range("M20").Value = 1.221658
$ws.Range("N20").Value = 2.443316
$ws.Range("O20").Value = 0.005117086949542552
$ws.Range("P20").Value = 0.003417220037046797
$ws.Range("Q20").Value = 1.034281928447
$ws.Range("R20").Value = 4.137127713788
$ws.Range("S20").Value = [double]"8.59838950898801E-05"
$ws.Range("T20").Value = [double]"3.849597882288967E-05"
$ws.Range("G21").Value = 0.8466215
$ws.Range("H21").Value = 1.693243
$ws.Range("I21").Value = 0.01680328982832053
$ws.Range("J21").Value = 0.0112652912032432
$ws.Range("M21").Value = 100.7253213333333
$ws.Range("N21").Value = 302.175964
$ws.Range("O21").Value = 0.4219022240948613
$ws.Range("P21").Value = 0.4226230904617871
$ws.Range("Q21").Value = 85.27622263520867
$ws.Range("R21").Value = 511.657335811252
$ws.Range("S21").Value = 0.007089345350678991
$ws.Range("T21").Value = 0.004760972183266623
$ws.Range("G22").Value = 4.273823333333334
$ws.Range("H22").Value = 12.82147
$ws.Range("I22").Value = 0.08482455518202518
$ws.Range("J22").Value = 0.08530234183968073
$ws.Range("M22").Value = 0.5804443333333333
$ws.Range("N22").Value = 1.741333
$ws.Range("O22").Value = 0.002431273010151717
$ws.Range("P22").Value = 0.002435427107574628
$ws.Range("Q22").Value = 2.480716535501111
$ws.Range("R22").Value = 22.32644881951
$ws.Range("S22").Value = 0.0002062316516121828
$ws.Range("T22").Value = 0.0002077476356559559
$ws.Range("G23").Value = 4.273823333333334
$ws.Range("H23").Value = 12.82147
$ws.Range("I23").Value = 0.08482455518202518
$ws.Range("J23").Value = 0.08530234183968073
$ws.Range("O23").Value = 0.0004752041289926495
$ws.Range("P23").Value = 0.00047601606752829
$ws.Range("Q23").Value = 0.4848681063822222
$ws.Range("R23").Value = 4.36381295744
$ws.Range("S23").Value = [double]"4.030897886246321E-05"
$ws.Range("T23").Value = [double]"4.060528531347874E-05"
$ws.Range("G24").Value = 4.273823333333334
$ws.Range("H24").Value = 12.82147
$ws.Range("I24").Value = 0.08482455518202518
$ws.Range("J24").Value = 0.08530234183968073
$ws.Range("M24").Value = 136.1000366666667
$ws.Range("N24").Value = 408.30011
$ws.Range("O24").Value = 0.5700742118164518
$ws.Range("P24").Value = 0.5710482463260632
$ws.Range("Q24").Value = 581.6675123735223
$ws.Range("R24").Value = 5235.0076113617
$ws.Range("S24").Value = 0.04835629143807413
$ws.Range("T24").Value = 0.04871175271505605
$ws.Range("G25").Value = 4.273823333333334
$ws.Range("H25").Value = 12.82147
$ws.Range("I25").Value = 0.08482455518202518
$ws.Range("J25").Value = 0.08530234183968073
$ws.Range("M25").Value = 1.221658
$ws.Range("N25").Value = 2.443316
$ws.Range("O25").Value = 0.005117086949542552
$ws.Range("P25").Value = 0.003417220037046797
$ws.Range("Q25").Value = 5.221150465753333
$ws.Range("R25").Value = 31.32690279452
$ws.Range("S25").Value = 0.0004340546243226931
$ws.Range("T25").Value = 0.0002914968717415723
$ws.Range("G26").Value = 4.273823333333334
$ws.Range("H26").Value = 12.82147
$ws.Range("I26").Value = 0.08482455518202518
$ws.Range("J26").Value = 0.08530234183968073
$ws.Range("M26").Value = 100.7253213333333
$ws.Range("N26").Value = 302.175964
$ws.Range("O26").Value = 0.4219022240948613
$ws.Range("P26").Value = 0.4226230904617871
$ws.Range("Q26").Value = 430.4822285718978
$ws.Range("R26").Value = 3874.34005714708
$ws.Range("S26").Value = 0.03578766848915371
$ws.Range("T26").Value = 0.06137193925321967
